$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112, shifting existing rows 112-150 down to 113-151.
$ws.Rows.Item(112).Insert()

# Populate the constant columns (same values used across the whole dataset).
$ws.Cells.Item(112, 1).Value = 7
$ws.Cells.Item(112, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(112, 3).Value = "Ñuble"
$ws.Cells.Item(112, 5).Value = 16
$ws.Cells.Item(112, 6).Value = 100112043
$ws.Cells.Item(112, 7).Value = "Pepino ensalada"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 18).Value = "Hortaliza"

# Populate the new row's specific data.
$ws.Cells.Item(112, 4).Value = 44468
$ws.Cells.Item(112, 10).Value = 120
$ws.Cells.Item(112, 11).Value = 16000
$ws.Cells.Item(112, 12).Value = 17000
$ws.Cells.Item(112, 13).Value = 16500
$ws.Cells.Item(112, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(112, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(112, 16).Value = 275
$ws.Cells.Item(112, 17).Value = 60
